# Updated BOM to source TLV431 from Newark.
#
# The TLV431 shunt-regulator line (row 6) previously sourced the
# "TLV431CDBVR" part directly from TI with the TI part number duplicated
# into the Vendor / Vendor part columns. Re-point it at Newark, using
# Newark's own vendor part number and the "AI" graded manufacturer part,
# and update the unit cost accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel's reference style toggled to R1C1 (cosmetic workbook calc setting).
$excel.ReferenceStyle = -4150   # xlR1C1

# Row 6: TLV431 shunt regulator - now sourced from Newark.
# Write G6 before E6 so new shared-string entries are appended in the
# same order Excel produced them (Vendor part, then Mfg part).
$ws.Range("G6").Value = "76C8858"        # Vendor part (Newark)
$ws.Range("E6").Value = "TLV431AIDBVR"   # Mfg part
$ws.Range("F6").Value = "Newark"         # Vendor
$ws.Range("H6").Value = 0.328            # Unit cost
$ws.Range("E6").Style = "Normal"         # clear the leftover text-format style

# Selection moved from I7 to H7.
$ws.Range("H7").Select()
